# B6-PowerPoint.pptx edit:
#  1. Re-style the three DrawingML tables (slides 14, 15, 16) from the
#     deck's default "Table_0" style to the
#     {2452237C-231C-48D4-9319-7898F90D55E6} table style.
#  2. Swap the presentation's applied theme colours from the "Integral"
#     (Red Violet) palette to the standard "Office Theme" palette - this
#     is how the underlying theme part content was swapped.

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$newTableStyle = "{2452237C-231C-48D4-9319-7898F90D55E6}"

foreach ($slideIdx in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIdx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable -eq -1) {
            $shp.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2. Theme colours --------------------------------------------------
# Office Theme colour scheme values (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink), replacing the Integral/"Red Violet" scheme.
$officeThemeRGB = @(
    0x000000,   # dk1
    0xFFFFFF,   # lt1
    0x44546A,   # dk2
    0xE7E6E6,   # lt2
    0x5B9BD5,   # accent1
    0xED7D31,   # accent2
    0xA5A5A5,   # accent3
    0xFFC000,   # accent4
    0x4472C4,   # accent5
    0x70AD47,   # accent6
    0x0563C1,   # hlink
    0x954F72    # folHlink
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $themeColors.Count; $i++) {
    $hex = $officeThemeRGB[$i - 1]
    $r = [math]::Floor($hex / 0x10000) % 0x100
    $g = [math]::Floor($hex / 0x100) % 0x100
    $b = $hex % 0x100
    $themeColors.Item($i).RGB = ($b * 0x10000) + ($g * 0x100) + $r
}
